# edit.ps1 - apply the zhengzhong_reading_group.pptx update:
#  1. Update the cached "datetimeFigureOut" footer-date field from
#     10/19/2018 -> 10/30/2018 on the slide master and all 11 slide
#     layouts.
#  2. On the last slide (the "ruler of LSTM LM" slide), move/resize the
#     bullet textbox and append "(Section N)" citations to each bullet.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

# ---------------------------------------------------------------------
# 1. Date placeholder text: 10/19/2018 -> 10/30/2018
# ---------------------------------------------------------------------

# Slide master (the date placeholder is shape #3 there).
$m.Shapes.Item(3).TextFrame.TextRange.Text = "10/30/2018"

# Each slide layout's date placeholder shape index (found by locating the
# shape whose placeholder type is "dt" in each layout).
$dateShapeIndexByLayout = @{
    1  = 3
    2  = 3
    3  = 3
    4  = 4
    5  = 6
    6  = 2
    7  = 1
    8  = 4
    9  = 4
    10 = 3
    11 = 3
}

for ($li = 1; $li -le $m.CustomLayouts.Count; $li++) {
    $layout = $m.CustomLayouts.Item($li)
    $shapeIdx = $dateShapeIndexByLayout[$li]
    $layout.Shapes.Item($shapeIdx).TextFrame.TextRange.Text = "10/30/2018"
}

# ---------------------------------------------------------------------
# 2. Last slide: move/resize the bullet textbox + append citations
# ---------------------------------------------------------------------

$s = $p.Slides.Item($p.Slides.Count)
$sh = $s.Shapes.Item(2)

# Reposition / resize (height is unchanged - spAutoFit keeps it put).
$sh.Left = 345.8692913386
$sh.Top = 75.2268503937
$sh.Width = 429.6735

$tr = $sh.TextFrame.TextRange

$shift = 0

# Bullet 1
$c = $tr.Characters(1 + $shift, 78)
$c.Text = "Context before 20: order of tokens are important (shuffling causes huge loss). (Section 5.1)"
$shift = $shift + (92 - 78)

# Bullet 2
$c = $tr.Characters(81 + $shift, 89)
$c.Text = "Context before 50: LSTM without cache may copy word from such context to generate target. (Section and 6.1)"
$shift = $shift + (107 - 89)

# Bullet 3
$c = $tr.Characters(172 + $shift, 104)
$c.Text = "Context from 20 to 200: order of tokens is not important, but the show up of relevant words is important. (Section 5.1)"
$shift = $shift + (119 - 104)

# Bullet 4
$c = $tr.Characters(278 + $shift, 74)
$c.Text = "Context before ~200: effective size of LSTM language model without cache. (Section 4)"
$shift = $shift + (85 - 74)

# Bullet 5
$c = $tr.Characters(354 + $shift, 68)
$c.Text = "Context from 50 to 300: only a rough semantic context is maintained. (Section 6.1)"
$shift = $shift + (82 - 68)

# Bullet 6
$c = $tr.Characters(424 + $shift, 103)
$c.Text = "Context between 200~3750: cache can help LSTM language model to retrieve information from the history. (Section 6.2)"
$shift = $shift + (116 - 103)
